$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update header row text:
# "SalesID" -> "SalesOrder"
# "SemlineNumber" -> "SalelineNumber"
$ws.Range("C1").Value = "SalesOrder"
$ws.Range("D1").Value = "SalelineNumber"
